$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 66
$ws.Range("G6").Value = 1972.08
$ws.Range("B10").Value = 27233.35
$ws.Range("F21").Value = 147
$ws.Range("G21").Value = 3773.49
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 1136.64
$ws.Range("F28").Value = 68
$ws.Range("G28").Value = 3032.8
$ws.Range("B32").Value = 12425.95
$ws.Range("F70").Value = 7
$ws.Range("G70").Value = 944.65
$ws.Range("F71").Value = 316
$ws.Range("G71").Value = 20129.2
$ws.Range("F78").Value = 39
$ws.Range("G78").Value = 2219.1
$ws.Range("F86").Value = 57
$ws.Range("G86").Value = 7151.79
$ws.Range("B90").Value = 173641.77
$ws.Range("F115").Value = 193
$ws.Range("G115").Value = 18684.33
$ws.Range("B117").Value = 12542.21
$ws.Range("B127").Value = 57552
$ws.Range("E127").Value = 136.86
$ws.Range("F127").Value = -5
$ws.Range("G127").Value = -603.45
$ws.Range("B128").Value = 64329
$ws.Range("E128").Value = 128.32
$ws.Range("F128").Value = 1
$ws.Range("G128").Value = 120.69
$ws.Range("F141").Value = 44
$ws.Range("G141").Value = 2355.32
$ws.Range("B142").Value = 2867
$ws.Range("F150").Value = 33
$ws.Range("G150").Value = 1534.17
$ws.Range("B156").Value = 31010.26
$ws.Range("F218").Value = 5
$ws.Range("G218").Value = 1081.1
$ws.Range("F222").Value = 11
$ws.Range("G222").Value = 1594.23
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 65
$ws.Range("G227").Value = 9378.200000000001
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("F229").Value = 55
$ws.Range("G229").Value = 7891.4
$ws.Range("F255").Value = 536
$ws.Range("G255").Value = 91832.88
$ws.Range("B260").Value = 177978.24
$ws.Range("F280").Value = 132
$ws.Range("G280").Value = 22326.48
$ws.Range("F288").Value = 38
$ws.Range("G288").Value = 3533.62
$ws.Range("F292").Value = 42
$ws.Range("G292").Value = 3497.34
$ws.Range("B304").Value = 170752.25
$ws.Range("F306").Value = 64
$ws.Range("G306").Value = 1351.04
$ws.Range("B309").Value = 1772.81
$ws.Range("F333").Value = 40
$ws.Range("G333").Value = 1995.6
$ws.Range("F334").Value = 191
$ws.Range("G334").Value = 9897.620000000001
$ws.Range("F338").Value = 76
$ws.Range("G338").Value = 1801.2
$ws.Range("F342").Value = 139
$ws.Range("G342").Value = 4402.13
$ws.Range("F345").Value = 45
$ws.Range("G345").Value = 2763.45
$ws.Range("B346").Value = 25222.2
$ws.Range("F354").Value = 13
$ws.Range("G354").Value = 891.67
$ws.Range("B358").Value = 34914.66
$ws.Range("F422").Value = 10
$ws.Range("G422").Value = 2161.1
$ws.Range("B424").Value = 2593.77
$ws.Range("F453").Value = 18
$ws.Range("G453").Value = 477.18
$ws.Range("F454").Value = 48
$ws.Range("G454").Value = 1639.2
$ws.Range("F456").Value = 43
$ws.Range("G456").Value = 4754.51
$ws.Range("B460").Value = 12725.18
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 107
$ws.Range("G473").Value = 3512.81
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F484").Value = 0
$ws.Range("G484").Value = 0
$ws.Range("F485").Value = 11
$ws.Range("G485").Value = 1930.17
$ws.Range("B488").Value = 29455.4
$ws.Range("F509").Value = 207
$ws.Range("G509").Value = 16638.66
$ws.Range("B510").Value = 22355.36
$ws.Range("F549").Value = 23
$ws.Range("G549").Value = 1100.78
$ws.Range("B560").Value = 4309.89
$ws.Range("F577").Value = 45
$ws.Range("G577").Value = 1934.55
$ws.Range("F579").Value = 33
$ws.Range("G579").Value = 2659.8
$ws.Range("F580").Value = 50
$ws.Range("G580").Value = 2849.5
$ws.Range("B583").Value = 14072.18
$ws.Range("F599").Value = 1480
$ws.Range("G599").Value = 241402.8
$ws.Range("F601").Value = 382
$ws.Range("G601").Value = 108056.34
$ws.Range("B606").Value = 397318.44
$ws.Range("B619").Value = 1664557.96
$ws.Range("B620").Value = 1664557.96
